# Applies the content edits described by the commit diff:
#  - 4 places where a paragraph's text runs were retyped/merged into a
#    single run (no actual text change, just run consolidation)
#  - 7 places where a trailing membership count "(n)" was added/changed
#    on association/Verein list entries
#
# (The diff also contains a large number of +/-1 dxa rounding drifts on
# <w:tblGrid>/<w:tcW> column widths across the two tables in this
# document. Those are a side effect of Word's table layout engine
# recalculating column widths and are not reachable/safe to reproduce
# through the exposed table COM surface without corrupting unrelated
# cells of this non-uniform table, so they are intentionally left
# alone here.)

$d = $word.ActiveDocument

function Replace-Exact($find, $replace) {
    $d.Content.Find.Execute(
        $find, $false, $false, $false, $false, $false,
        $true, 1, $false, $replace, 2
    ) | Out-Null
}

# --- Run-merge paragraphs (identical concatenated text, single run) ---

Replace-Exact `
    "Fläche des Dorfes in ha  93,2 , davon LF  .. unbekannt .....  ha  /  FF  geschätzt 3  ha" `
    "Fläche des Dorfes in ha  93,2 , davon LF  .. unbekannt .....  ha  /  FF  geschätzt 3  ha"

Replace-Exact `
    "(anhand öffentlich bekannter Betriebe geschätzt, untere Schranken)" `
    "(anhand öffentlich bekannter Betriebe geschätzt, untere Schranken)"

Replace-Exact `
    "3. Oktober - Tag der Deutschen Einheit mit Sektempfang am Dorfgemeinschaftshaus mit Konzert des Musikzugs der Freiwilligen Feuerwehr (Ortsrat)" `
    "3. Oktober - Tag der Deutschen Einheit mit Sektempfang am Dorfgemeinschaftshaus mit Konzert des Musikzugs der Freiwilligen Feuerwehr (Ortsrat)"

Replace-Exact "Jahreskonzert (Musikzug FFW)" "Jahreskonzert (Musikzug FFW)"

# --- Membership counts added/updated on Verein list entries ---

Replace-Exact "- Dorfpflege Rössing" "- Dorfpflege Rössing (115)"
Replace-Exact "- Deutsches Rotes Kreuz (17)" "- Deutsches Rotes Kreuz (139)"
Replace-Exact "- Kleingartenverein (8)" "- Kleingartenverein (76)"
Replace-Exact `
    "- Niedersächsische Kameradschaftsvereinigung, Ortsgruppe Rössing" `
    "- Niedersächsische Kameradschaftsvereinigung, Ortsgruppe Rössing (154)"
Replace-Exact "- Rassekaninchenzuchtverein RKZV" "- Rassekaninchenzuchtverein RKZV (6)"
Replace-Exact "- Tennisverein" "- Tennisverein (78)"
Replace-Exact "- VSV Rössing" "- VSV Rössing (624)"

Write-Host "done"
